$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "MLX90614:" section, mirrors the existing "SHT31:" section (rows 6-10) ---

# Row 12: section header + instructions (matches row 6 formatting)
$ws.Range("A12").Value = "MLX90614:"
$ws.Range("A12").Font.Bold = $true
$ws.Range("B12").Value = $ws.Range("B6").Value()

# Row 13: VIN / Solder / Red Sub-Wire / Screw / "3.3V" on RHS of LEMS Shield
$ws.Range("A13").Value = $ws.Range("A7").Value()
$ws.Range("B13").Value = $ws.Range("B7").Value()
$ws.Range("C13").Value = $ws.Range("C7").Value()
$ws.Range("D13").Value = $ws.Range("D7").Value()
$ws.Range("E13").Value = "'" + $ws.Range("E7").Value()

# Row 14: GND / Solder / Black Sub-Wire / Screw / "GND" on RHS of LEMS Shield
$ws.Range("A14").Value = $ws.Range("A8").Value()
$ws.Range("B14").Value = $ws.Range("B8").Value()
$ws.Range("C14").Value = $ws.Range("C8").Value()
$ws.Range("D14").Value = $ws.Range("D8").Value()
$ws.Range("E14").Value = "'" + $ws.Range("E8").Value()

# Row 15: SCL / Solder / White Sub-Wire / Screw / "SCL" on RHS of LEMS Shield
$ws.Range("A15").Value = $ws.Range("A9").Value()
$ws.Range("B15").Value = $ws.Range("B9").Value()
$ws.Range("C15").Value = $ws.Range("C9").Value()
$ws.Range("D15").Value = $ws.Range("D9").Value()
$ws.Range("E15").Value = $ws.Range("E9").Value()

# Row 16: SDA / Solder / Green Sub-Wire / Screw / "SDA" on RHS of LEMS Shield
$ws.Range("A16").Value = $ws.Range("A10").Value()
$ws.Range("B16").Value = $ws.Range("B10").Value()
$ws.Range("C16").Value = $ws.Range("C10").Value()
$ws.Range("D16").Value = $ws.Range("D10").Value()
$ws.Range("E16").Value = $ws.Range("E10").Value()

# Keep the active-cell selection in sync with where Excel would land after the edit
$ws.Range("E18").Select()
